$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number + report week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  35"
$ws.Range("C9").Value = "Report Covering the Week  8/28/2023  Through  9/3/2023"

# --- Weekly crime-stat table updates (rows 14-29) ---
$ws.Range("M14").Value = -68.421052631578
$ws.Range("N14").Value = -87.5
$ws.Range("F15").Value = 5
$ws.Range("H15").Value = 66.666666666666
$ws.Range("I15").Value = 29
$ws.Range("J15").Value = 22
$ws.Range("K15").Value = 31.818181818181
$ws.Range("L15").Value = 3.571428571428
$ws.Range("M15").Value = 45
$ws.Range("N15").Value = -53.225806451612
$ws.Range("C16").Value = 2
$ws.Range("E16").Value = -60
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = -56.25
$ws.Range("I16").Value = 179
$ws.Range("J16").Value = 241
$ws.Range("K16").Value = -25.726141078838
$ws.Range("L16").Value = 25.174825174825
$ws.Range("M16").Value = -40.531561461794
$ws.Range("N16").Value = -88.262295081967
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 57
$ws.Range("H17").Value = -21.052631578947
$ws.Range("I17").Value = 488
$ws.Range("J17").Value = 488
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 17.026378896882
$ws.Range("M17").Value = 12.962962962963
$ws.Range("N17").Value = -39.152119700748
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 25
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 83
$ws.Range("J18").Value = 132
$ws.Range("K18").Value = -37.121212121212
$ws.Range("L18").Value = -7.777777777777
$ws.Range("M18").Value = -57.868020304568
$ws.Range("N18").Value = -84.514925373134
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 230
$ws.Range("J19").Value = 270
$ws.Range("K19").Value = -14.814814814814
$ws.Range("L19").Value = 19.791666666666
$ws.Range("M19").Value = -21.501706484641
$ws.Range("N19").Value = -63.549920760697
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -38.461538461538
$ws.Range("I20").Value = 107
$ws.Range("J20").Value = 129
$ws.Range("K20").Value = -17.054263565891
$ws.Range("L20").Value = 10.309278350515
$ws.Range("M20").Value = 0.943396226415
$ws.Range("N20").Value = -74.764150943396
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 34
$ws.Range("E21").Value = -26.470588235294
$ws.Range("F21").Value = 99
$ws.Range("G21").Value = 139
$ws.Range("H21").Value = -28.776978417266
$ws.Range("I21").Value = 1122
$ws.Range("J21").Value = 1295
$ws.Range("K21").Value = -13.359073359073
$ws.Range("L21").Value = 14.489795918367
$ws.Range("M21").Value = -17.982456140350
$ws.Range("N21").Value = -72.144985104270
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 4
$ws.Range("I22").Value = 32
$ws.Range("J22").Value = 37
$ws.Range("K22").Value = -13.513513513513
$ws.Range("L22").Value = 23.076923076923
$ws.Range("M22").Value = -3.030303030303
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -66.666666666666
$ws.Range("F23").Value = 24
$ws.Range("G23").Value = 26
$ws.Range("H23").Value = -7.692307692307
$ws.Range("I23").Value = 243
$ws.Range("J23").Value = 265
$ws.Range("K23").Value = -8.301886792452
$ws.Range("L23").Value = 14.084507042253
$ws.Range("M23").Value = 48.170731707317
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -11.764705882352
$ws.Range("F24").Value = 75
$ws.Range("G24").Value = 79
$ws.Range("H24").Value = -5.063291139240
$ws.Range("I24").Value = 765
$ws.Range("J24").Value = 781
$ws.Range("K24").Value = -2.048655569782
$ws.Range("L24").Value = 47.115384615384
$ws.Range("M24").Value = 14.864864864864
$ws.Range("C25").Value = 27
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = 28.571428571428
$ws.Range("F25").Value = 92
$ws.Range("G25").Value = 71
$ws.Range("H25").Value = 29.577464788732
$ws.Range("I25").Value = 669
$ws.Range("J25").Value = 654
$ws.Range("K25").Value = 2.293577981651
$ws.Range("L25").Value = 24.581005586592
$ws.Range("M25").Value = -27.440347071583
$ws.Range("F26").Value = 6
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 50
$ws.Range("I26").Value = 34
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = -5.555555555555
$ws.Range("L26").Value = -26.086956521739
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 125
$ws.Range("I27").Value = 62
$ws.Range("K27").Value = 6.896551724137
$ws.Range("L27").Value = 6.896551724137
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = -46.666666666666
$ws.Range("L28").Value = -44.827586206896
$ws.Range("M28").Value = -60
$ws.Range("N28").Value = -85.777777777777
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = -75
$ws.Range("J29").Value = 47
$ws.Range("K29").Value = -36.170212765957
$ws.Range("L29").Value = -36.170212765957
$ws.Range("M29").Value = -53.846153846153
$ws.Range("N29").Value = -85.576923076923

# --- Cells whose type/style changes (number <-> "N/A" text) ---
$ws.Range("F14").Value = "'0"
$ws.Range("N22").Copy()
$ws.Range("F14").PasteSpecial(-4122)
$ws.Range("C15").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = 100
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").Value = 0
$ws.Range("K14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("G22").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = 300
$ws.Range("K14").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("C26").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 100
$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("C27").Value = "'0"
$ws.Range("N22").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$ws.Range("N22").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = "'***.*"
$ws.Range("N22").Copy()
$ws.Range("E27").PasteSpecial(-4122)
